# Auto-generated PowerShell/COM script applying the cryptos.xlsx diff
# (crypto price/volume refresh commit, GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.452.21'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '2.489.64'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('E4').Value = '  -0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '569.00'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.43%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '164.17'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -1.18%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.510'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('D9').Value = '2.487.24'
$ws.Range('E9').Value = '  -0.80%  '
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('E11').Value = '  -0.60%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.352'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('D14').Value = '2.944.16'
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('D15').Value = '69.242.21'
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('E16').Value = '  -0.26%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '24.08'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -2.97%  '
$ws.Range('D18').Value = '2.491.96'
$ws.Range('E18').Value = '  -0.92%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '11.14'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -1.78%  '
$ws.Range('E20').Value = '  -3.74%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '346.49'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('E22').Value = '  -1.25%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '1.88'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -3.98%  '
$ws.Range('E24').Value = '  -0.07%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '69.22'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -1.08%  '
$ws.Range('E26').Value = '  -2.62%  '
$ws.Range('D27').Value = '2.616.20'
$ws.Range('E27').Value = '  -1.21%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '8.59'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -3.55%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('D30').Value = '0.0₃0864'
$ws.Range('E30').Value = '  -3.08%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '7.53'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -4.04%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '437.88'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -5.18%  '
$ws.Range('E33').Value = '  -4.38%  '
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('E35').Value = '  -1.59%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '156.17'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('E37').Value = '  -3.28%  '
$ws.Range('E38').Value = '  +0.46%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '18.10'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -2.24%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  -1.94%  '
$ws.Range('E42').Value = '  -3.87%  '
$ws.Range('E43').Value = '  -1.85%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '2.33'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +46.50%  '
$ws.Range('E45').Value = '  -5.03%  '
$ws.Range('E46').Value = '  -6.13%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '137.85'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -2.57%  '
$ws.Range('E48').Value = '  -1.98%  '
$ws.Range('E49').Value = '  -4.19%  '
$ws.Range('E50').Value = '  -0.68%  '
$ws.Range('E51').Value = '  -0.83%  '
